# DriveController & UART for CP2102
# - New IMU-adjacent peripheral pin mapping for the CP2102 UART (PB02/PB03)
# - Shares the W25Q128 SPI bus pins with the newly added ADS1118 SPI device

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- W25Q128.SPI pins are now shared with ADS1118.SPI -----------------------
# Update the 4 cells that previously only said "W25Q128.SPI" so they reflect
# both devices sharing the bus. Apply the smaller wrap-capable font so the
# two-line label fits, and grow the row height to match.
$newSpiLabel = "W25Q128.SPI" + [char]10 + "ADS1118.SPI"

# Fully format the first cell, then propagate the exact same style to the
# other three via copy/paste-special so the style table doesn't accumulate
# extra transient combinations (one per cell instead of one shared union).
$ws.Range("R11").Value = $newSpiLabel
$ws.Range("R11").Font.Size = 8
$ws.Range("R11").Interior.Color = 49407
$ws.Range("R11").WrapText = $true

$ws.Range("R11").Copy()
$ws.Range("O12").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("O14").PasteSpecial(-4122)
$ws.Range("R14").PasteSpecial(-4122)

$ws.Range("O12").Value = $newSpiLabel
$ws.Range("O14").Value = $newSpiLabel
$ws.Range("R14").Value = $newSpiLabel

$ws.Rows.Item(11).RowHeight = 23.25
$ws.Rows.Item(12).RowHeight = 23.25
$ws.Rows.Item(14).RowHeight = 23.25

# --- New CP2102.UART peripheral on PB02/PB03 (O6/R6) -------------------------
$ws.Range("O6").Value = "CP2102.UART"
$ws.Range("O6").Interior.Color = 15773696
$ws.Range("R6").Value = "CP2102.UART"
$ws.Range("R6").Interior.Color = 15773696

# --- Restore cursor/selection -------------------------------------------
[void]$ws.Range("N21").Select()
